$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (fcs) values for rows 2-7
$ws.Range("C2").Value = 12626.2716881561
$ws.Range("C3").Value = 11721.816447481
$ws.Range("C4").Value = 7791.75184438609
$ws.Range("C5").Value = 7816.4328508422
$ws.Range("C6").Value = 11813.199742613
$ws.Range("C7").Value = 12111.5808042358

# Update column F (need_to_buy) values for rows 2-7
$ws.Range("F2").Value = -5.06321580757761
$ws.Range("F3").Value = 346.230448992106
$ws.Range("F4").Value = 155.50459288658
$ws.Range("F5").Value = 151.778700628926
$ws.Range("F6").Value = 324.787171251434
$ws.Range("F7").Value = 347.919045381132

# Update column E (MYDIR) values for rows 12-15
$ws.Range("E12").Value = 7892.84507689171
$ws.Range("E13").Value = 7808.44073014521
$ws.Range("E14").Value = 7808.44073014521
$ws.Range("E15").Value = 7808.44073014521

# Update column F (need_to_buy) values for rows 12-15
$ws.Range("F12").Value = 246.901317391531
$ws.Range("F13").Value = 382.065785276785
$ws.Range("F14").Value = 379.975636314554
$ws.Range("F15").Value = 374.864269370142
